$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 146; this shifts the existing rows 146-186 down to 147-187.
$ws.Rows.Item(146).Insert()

# Populate the newly inserted row 146 with the new weekly record.
$ws.Range("A146").Value = 4
$ws.Range("B146").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C146").Value = "Los Lagos"
$ws.Range("D146").Value = 44588
$ws.Range("E146").Value = 10
$ws.Range("F146").Value = 100112032
$ws.Range("G146").Value = "Zapallo italiano"
$ws.Range("H146").Value = "Sin especificar"
$ws.Range("I146").Value = "Primera"
$ws.Range("J146").Value = 140
$ws.Range("K146").Value = 15000
$ws.Range("L146").Value = 16000
$ws.Range("M146").Value = 15500
$ws.Range("N146").Value = "$/caja 50 unidades"
$ws.Range("O146").Value = "Región Metropolitana"
$ws.Range("P146").Value = 310
$ws.Range("Q146").Value = 50
$ws.Range("R146").Value = "Hortaliza"
